$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the "[onshow..now;frm=’yyyy-mm-dd hh:nn:ss’]" run into five runs:
#      [onshow..now;frm=  |  '  |  yyyy-mm-dd hh:nn:ss  |  '  |  ]
#    The curly single-quotes in the original become straight apostrophes,
#    each becoming its own run (all with identical rPr), matching how Word
#    leaves behind separate same-format runs when the text is retyped.
# ---------------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("[onshow..now;frm=’yyyy-mm-dd hh:nn:ss’]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Replace the whole run's text (keeps the run's existing rPr).
    $rng.Text = "[onshow..now;frm='yyyy-mm-dd hh:nn:ss']"

    $base = $rng.Start

    # Offsets of the two apostrophes and the quoted format string, within
    # the replaced text "[onshow..now;frm='yyyy-mm-dd hh:nn:ss']".
    $q1Start = $base + 17
    $q1End   = $base + 18
    $fmtStart = $base + 18
    $fmtEnd   = $base + 37
    $q2Start = $base + 37
    $q2End   = $base + 38

    # Toggling a character property on/off on each sub-range forces Word to
    # split it into its own run while leaving the final formatting (rPr)
    # identical to its neighbours.
    $p1 = $d.Range($q1Start, $q1End)
    $p1.Bold = 1
    $p1.Bold = 0

    $p2 = $d.Range($fmtStart, $fmtEnd)
    $p2.Bold = 1
    $p2.Bold = 0

    $p3 = $d.Range($q2Start, $q2End)
    $p3.Bold = 1
    $p3.Bold = 0
}

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from right after the chart drawing to right
#    after the "[onshow..cst.PHP_VERSION]" run.
# ---------------------------------------------------------------------------

$phpRng = $d.Content
$phpFound = $phpRng.Find.Execute("[onshow..cst.PHP_VERSION]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($phpFound) {
    $gbTarget = $d.Range($phpRng.End, $phpRng.End)
    $d.Bookmarks.Add("_GoBack", $gbTarget)
}
